$d = $word.ActiveDocument

$pairs = @(
    @("76×65=4940", "58×78=4524"),
    @("34×38=1292", "58×81=4698"),
    @("62×88=5456", "81×25=2025"),
    @("89×38=3382", "29×70=2030"),
    @("89×52=4628", "85×23=1955"),
    @("56×49=2744", "79×99=7821"),
    @("66×19=1254", "69×63=4347"),
    @("22×82=1804", "90×13=1170"),
    @("92×26=2392", "32×31=992"),
    @("35×43=1505", "44×70=3080"),
    @("97×44=4268", "37×33=1221"),
    @("36×43=1548", "93×19=1767"),
    @("67×45=3015", "39×20=780"),
    @("39×46=1794", "42×46=1932"),
    @("15×68=1020", "90×14=1260"),
    @("72×94=6768", "92×51=4692"),
    @("46×66=3036", "97×23=2231"),
    @("75×42=3150", "55×85=4675"),
    @("29×29=841", "47×60=2820"),
    @("11×18=198", "82×68=5576"),
    @("51×38=1938", "21×95=1995"),
    @("68×14=952", "74×12=888"),
    @("75×72=5400", "18×90=1620"),
    @("50×85=4250", "62×47=2914"),
    @("81×40=3240", "28×86=2408")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
